$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fix: Agregar cohorte para produccion de ARL"
# Column K ("mes_cotizacion") held a plain yyyymm integer (202412) as a placeholder.
# Replace it with the real cohort date (2024-12-31), formatted the same way as the
# other date column (J) in this sheet: yyyy-mm-dd.
$ws.Range("K2:K4").Value = 45657
$ws.Range("K2:K4").NumberFormat = "yyyy\-mm\-dd;@"

# Reflect where the author was last looking in the sheet: scrolled over so column C
# is the left-most visible column, with cell K4 selected/active.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = $ws.Range("C1").Column
$ws.Range("K4").Select()
